$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").ClearContents()
$ws.Range("B1").Value = "LED"
$ws.Range("C1").Value = "Buzzer"
$ws.Range("D1").Value = "blinking"

# Row 2 - xmin
$ws.Range("A2").Value = "xmin"
$ws.Range("B2").Value = 3.5
$ws.Range("C2").Formula = "=B2"
$ws.Range("D2").Formula = "=C2"

# Row 3 - xmax
$ws.Range("A3").Value = "xmax"
$ws.Range("B3").Value = 40
$ws.Range("C3").Formula = "=B3"
$ws.Range("D3").Formula = "=C3"

# Row 4 - ymin
$ws.Range("A4").Value = "ymin"
$ws.Range("B4").Value = 511
$ws.Range("C4").Value = 49
$ws.Range("D4").Value = 24462

# Row 5 - ymax
$ws.Range("A5").Value = "ymax"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = 489237

# Row 6 - m
$ws.Range("A6").Value = "m"
$ws.Range("B6").Formula = "=(B5-B4)/(B3-B2)"
$ws.Range("C6").Formula = "=(C5-C4)/(C3-C2)"
$ws.Range("D6").Formula = "=(D5-D4)/(D3-D2)"

# Row 7 - b
$ws.Range("A7").Value = "b"
$ws.Range("B7").Formula = "=B4-B2*B6"
$ws.Range("C7").Formula = "=C4-C2*C6"
$ws.Range("D7").Formula = "=D4-D2*D6"

# Row 8
$ws.Range("D8").Formula = "=LOG(D5,2)-1"

# Row 9 - fmax
$ws.Range("A9").Value = "fmax"
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = "kHz"
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = "Hz"

# Row 10 - fmin
$ws.Range("A10").Value = "fmin"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "kHz"
$ws.Range("F10").Value = 0.2
$ws.Range("G10").Value = "Hz"

# Row 11 - count
$ws.Range("A11").Value = "count"
$ws.Range("C11").Value = 511
$ws.Range("D11").Value = "clock cycles"
$ws.Range("F11").Formula = "=C11"
$ws.Range("G11").Formula = "=D11"

# Row 12 - clk freq
$ws.Range("A12").Value = "clk freq"
$ws.Range("C12").Value = 50000
$ws.Range("D12").Value = "kHz"
$ws.Range("F12").Formula = "=C12*1000"
$ws.Range("G12").Value = "Hz"

# Row 13 - PWM freq
$ws.Range("A13").Value = "PWM freq"
$ws.Range("C13").Formula = "=C12/C11"
$ws.Range("D13").Value = "kHz"
$ws.Range("F13").Formula = "=C13*1000"
$ws.Range("G13").Value = "Hz"

# Row 14 - fmax factor
$ws.Range("A14").Value = "fmax factor"
$ws.Range("C14").Formula = "=C13/C9"
$ws.Range("F14").Formula = "=F13/F9"

# Row 15 - fmin factor
$ws.Range("A15").Value = "fmin factor"
$ws.Range("C15").Formula = "=C13/C10"
$ws.Range("F15").Formula = "=F13/F10"

# Row 17
$ws.Range("C17").Value = 0.02
$ws.Range("D17").Value = "ms"

# Row 18
$ws.Range("B18").Value = "clk freq"
$ws.Range("C18").Formula = "=C12"
$ws.Range("D18").Formula = "=D12"
$ws.Range("F18").Value = "initial pwm period"
$ws.Range("H18").Value = 0.01024
$ws.Range("I18").Value = "ms"

# Row 19
$ws.Range("B19").Value = "clk period"
$ws.Range("C19").Formula = "=1/C18"
$ws.Range("D19").Value = "ms"
$ws.Range("F19").Value = "period"
$ws.Range("H19").Value = 3

# Row 20
$ws.Range("B20").Value = "period"
$ws.Range("C20").Formula = "=C17/C19"
$ws.Range("F20").Value = "expect pwm period"
$ws.Range("H20").Formula = "=H19*H18"

$ws.Range("I9").Select()
